# Update cryptocurrency price and volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.074.75"
$ws.Range("E2").Value = "  +4.02%  "

$ws.Range("D3").Value = "'2.649.28"
$ws.Range("E3").Value = "  +6.22%  "

$ws.Range("D5").Value = "'114.11"
$ws.Range("E5").Value = "  +8.02%  "

$ws.Range("D6").Value = "'326.54"
$ws.Range("E6").Value = "  +2.88%  "

$ws.Range("D7").Value = "'0.530"
$ws.Range("E7").Value = "  +2.33%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.558"
$ws.Range("E9").Value = "  +4.28%  "

$ws.Range("D10").Value = "'41.03"
$ws.Range("E10").Value = "  +5.97%  "

$ws.Range("D11").Value = "'20.19"
$ws.Range("E11").Value = "  -0.74%  "

$ws.Range("D12").Value = "'0.0823"
$ws.Range("E12").Value = "  +2.63%  "

$ws.Range("E13").Value = "  +0.78%  "

$ws.Range("D14").Value = "'7.40"
$ws.Range("E14").Value = "  +4.59%  "

$ws.Range("D15").Value = "'3.063.19"
$ws.Range("E15").Value = "  +6.16%  "

$ws.Range("D16").Value = "'2.645.77"
$ws.Range("E16").Value = "  +6.01%  "

$ws.Range("E17").Value = "  +5.41%  "

$ws.Range("D18").Value = "'49.966.48"
$ws.Range("E18").Value = "  +4.11%  "

$ws.Range("D19").Value = "'13.25"
$ws.Range("E19").Value = "  +3.20%  "

$ws.Range("D20").Value = "'6.80"
$ws.Range("E20").Value = "  +3.19%  "

$ws.Range("E21").Value = "  -2.11%  "

$ws.Range("E22").Value = "  +3.08%  "

$ws.Range("D23").Value = "'72.14"
$ws.Range("E23").Value = "  +1.59%  "

$ws.Range("D24").Value = "'276.47"
$ws.Range("E24").Value = "  +2.70%  "

$ws.Range("D25").Value = "'2.59"
$ws.Range("E25").Value = "  +3.16%  "

$ws.Range("D26").Value = "'26.86"
$ws.Range("E26").Value = "  +4.35%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").Value = "'10.06"
$ws.Range("E28").Value = "  +3.57%  "

$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").Value = "'36.29"
$ws.Range("E30").Value = "  +5.17%  "

$ws.Range("D31").Value = "'0.141"
$ws.Range("E31").Value = "  +1.57%  "

$ws.Range("D32").Value = "'50.33"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("E33").Value = "  +3.46%  "

$ws.Range("D34").Value = "'19.54"
$ws.Range("E34").Value = "  +2.87%  "

$ws.Range("D35").Value = "'0.0813"
$ws.Range("E35").Value = "  +5.39%  "

$ws.Range("E36").Value = "  +10.76%  "

$ws.Range("E37").Value = "  -0.09%  "

$ws.Range("E38").Value = "  +7.20%  "

$ws.Range("E39").Value = "  +8.32%  "

$ws.Range("E40").Value = "  +2.27%  "

$ws.Range("D41").Value = "'123.41"
$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").Value = "'22.05"
$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("D44").Value = "'0.0317"
$ws.Range("E44").Value = "  +5.18%  "

$ws.Range("D45").Value = "'2.086.27"
$ws.Range("E45").Value = "  +4.26%  "

$ws.Range("E46").Value = "  +6.45%  "

$ws.Range("D47").Value = "'2.30"
$ws.Range("E47").Value = "  +15.25%  "

$ws.Range("E48").Value = "  +5.52%  "

$ws.Range("D49").Value = "'9.16"
$ws.Range("E49").Value = "  +2.89%  "

$ws.Range("D50").Value = "'5.41"
$ws.Range("E50").Value = "  +4.95%  "

$ws.Range("D51").Value = "'59.90"
$ws.Range("E51").Value = "  +6.60%  "
